# Populate the "per_minute" (sheet2) and "per_poss" (sheet3) worksheets,
# which were previously empty placeholders, with per-36-minute and
# per-100-possession stat tables mirroring the layout already used on the
# "per_game" (sheet1) and "advanced" (sheet4) tabs. The "advanced" tab
# picks up two new columns, ORtg/DRtg, which per_poss also gains.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # per_game   (reference formatting source)
$ws2 = $wb.Worksheets.Item(2)   # per_minute
$ws3 = $wb.Worksheets.Item(3)   # per_poss

function Fill-Row {
    param($ws, $row, $startCol, $values)
    $col = $startCol
    foreach ($v in $values) {
        if ($v -ne $null) {
            $ws.Cells.Item($row, $col).Value = $v
        }
        $col = $col + 1
    }
}

# ---------------------------------------------------------------------
# 1) Bring over the bold/bordered/centered header style (style index 1
#    on the existing sheets) by copying formats from the fully built
#    per_game sheet, so new cells match the rest of the workbook.
#    Use a single source cell so PasteSpecial fills the whole target
#    range uniformly instead of tiling a multi-column source pattern.
# ---------------------------------------------------------------------
$ws1.Range("B1").Copy()
$ws2.Range("B1:AE1").PasteSpecial(-4122)
$ws3.Range("B1:AH1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2:A9").PasteSpecial(-4122)
$ws3.Range("A2:A9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) per_minute (sheet2) header + data, columns A:AE
# ---------------------------------------------------------------------
Fill-Row $ws2 1 2 @("Season", "Age", "Tm", "Lg", "Pos", "G", "GS", "MP", "FG", "FGA", "FG%", "3P", "3PA", "3P%", "2P", "2PA", "2P%", "FT", "FTA", "FT%", "ORB", "DRB", "TRB", "AST", "STL", "BLK", "TOV", "PF", "PTS", "RSPS")

for ($r = 2; $r -le 9; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 2
}

Fill-Row $ws2 2 2 @("2018-19", "23", "CHO", "NBA", "PG", 46, 3, 676, 3.9, 11.5, 0.343, 1.8, 6.4, 0.281, 2.1, 5.1, 0.421, 1.9, 2.4, 0.7609999999999999, 0.5, 2.8, 3.4, 6.4, 1.2, 0.1, 1.6, 2.6, 11.6, "RS")
Fill-Row $ws2 4 2 @("2019-20", "24", "CHO", "NBA", "PG", 63, 53, 2211, 6, 15.7, 0.382, 3.5, 9.5, 0.373, 2.4, 6.2, 0.397, 3.1, 3.8, 0.82, 0.7, 2.8, 3.5, 7.7, 1, 0.2, 2.9, 1.9, 18.6, "RS")
Fill-Row $ws2 6 2 @("2020-21", "25", "CHO", "NBA", "PG", 51, 44, 1531, 5.4, 14.5, 0.374, 3.8, 10.3, 0.373, 1.6, 4.3, 0.376, 2.7, 3.3, 0.821, 0.4, 2.8, 3.2, 6.5, 1.1, 0.1, 1.9, 1.7, 17.4, "RS")
Fill-Row $ws2 8 2 @("Career", $null, $null, "NBA", $null, 160, 100, 4418, 5.5, 14.6, 0.375, 3.4, 9.3, 0.363, 2.1, 5.3, 0.394, 2.8, 3.4, 0.8140000000000001, 0.6, 2.8, 3.4, 7.1, 1.1, 0.2, 2.4, 1.9, 17.1, "RS")

# ---------------------------------------------------------------------
# 3) per_poss (sheet3) header + data, columns A:AH
#    (adds ORtg in AF and DRtg in AG ahead of the trailing RSPS column)
# ---------------------------------------------------------------------
Fill-Row $ws3 1 2 @("Season", "Age", "Tm", "Lg", "Pos", "G", "GS", "MP", "FG", "FGA", "FG%", "3P", "3PA", "3P%", "2P", "2PA", "2P%", "FT", "FTA", "FT%", "ORB", "DRB", "TRB", "AST", "STL", "BLK", "TOV", "PF", "PTS", $null, "ORtg", "DRtg", "RSPS")

for ($r = 2; $r -le 9; $r++) {
    $ws3.Cells.Item($r, 1).Value = $r - 2
}

Fill-Row $ws3 2 2 @("2018-19", "23", "CHO", "NBA", "PG", 46, 3, 676, 5.3, 15.5, 0.343, 2.4, 8.7, 0.281, 2.9, 6.8, 0.421, 2.5, 3.3, 0.7609999999999999, 0.7, 3.8, 4.5, 8.7, 1.7, 0.1, 2.2, 3.5, 15.6, $null, 104, 115, "RS")
Fill-Row $ws3 4 2 @("2019-20", "24", "CHO", "NBA", "PG", 63, 53, 2211, 8.300000000000001, 21.8, 0.382, 4.9, 13.3, 0.373, 3.4, 8.6, 0.397, 4.3, 5.3, 0.82, 1, 3.9, 4.8, 10.7, 1.4, 0.3, 4.1, 2.6, 25.9, $null, 109, 116, "RS")
Fill-Row $ws3 6 2 @("2020-21", "25", "CHO", "NBA", "PG", 51, 44, 1531, 7.4, 19.7, 0.374, 5.2, 13.9, 0.373, 2.2, 5.8, 0.376, 3.7, 4.5, 0.821, 0.6, 3.8, 4.4, 8.800000000000001, 1.5, 0.2, 2.6, 2.3, 23.6, $null, 114, 116, "RS")
Fill-Row $ws3 8 2 @("Career", $null, $null, "NBA", $null, 160, 100, 4418, 7.5, 20.1, 0.375, 4.6, 12.8, 0.363, 2.9, 7.3, 0.394, 3.8, 4.7, 0.8140000000000001, 0.8, 3.8, 4.6, 9.699999999999999, 1.5, 0.3, 3.3, 2.6, 23.5, $null, 110, 116, "RS")

# ---------------------------------------------------------------------
# 4) advanced (sheet4): no value changes -- the two new shared strings
#    (ORtg/DRtg) added above are enough for the existing header cells
#    there to keep referring to the same logical text (PER, TS%, ...).
# ---------------------------------------------------------------------

Write-Host "per_minute and per_poss populated"
